$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "93.175.40"
$ws.Range("E2").Value = "  -2.97%  "
$ws.Range("D3").Value = "3.349.11"
$ws.Range("E3").Value = "  -2.31%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.24%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "231.76"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.98%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "620.11"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.67%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.38"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.69%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.385"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -5.05%  "
$ws.Range("E9").Value = "  -0.13%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.931"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.82%  "
$ws.Range("D11").Value = "3.354.85"
$ws.Range("E11").Value = "  -1.81%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "42.20"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.52%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.192"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.72%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.00"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.11%  "
$ws.Range("D15").Value = "92.994.14"
$ws.Range("E15").Value = "  -3.32%  "
$ws.Range("D16").Value = "3.969.44"
$ws.Range("E16").Value = "  -2.91%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000243"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.71%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.99"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -6.64%  "
$ws.Range("D19").Value = "3.346.28"
$ws.Range("E19").Value = "  -2.82%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.36"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.66%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.17"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.01%  "
$ws.Range("B22").Value = "SuiNetwork"
$ws.Range("C22").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.32"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.05%  "
$ws.Range("B23").Value = "BitcoinCash"
$ws.Range("C23").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "490.31"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.87%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.455"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -7.80%  "
$ws.Range("E25").Value = "  -6.01%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.10"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -5.68%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "89.64"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.65%  "
$ws.Range("D28").Value = "3.517.93"
$ws.Range("E28").Value = "  -3.18%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "11.61"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -5.14%  "
$ws.Range("E30").Value = "  -0.17%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "11.19"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -5.48%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.66"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.42%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.136"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.57%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.998"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.06%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.172"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -6.04%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "28.39"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -6.63%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.529"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -7.73%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "526.89"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.86%  "
$ws.Range("E39").Value = "  +0.03%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.35"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.86%  "
$ws.Range("E41").Value = "  -6.20%  "
$ws.Range("E42").Value = "  -2.70%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.881"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.33%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "24.03"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.84%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.68"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.38%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.56"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.16%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0403"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.88%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.38"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.55%  "
$ws.Range("E49").Value = "  -2.49%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "52.44"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.61%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.90"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.81%  "
